$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', RobustScaler()), ('selector', None),`n                ('model', SVC(C=5, kernel='poly', random_state=42))])"
$ws.Range("B2").Value = 0.7062745098039216
$ws.Range("C2").Value = "{'selector': None, 'scaler': RobustScaler(), 'model__kernel': 'poly', 'model__class_weight': None, 'model__C': 5}"
$ws.Range("D2").Value = 0.7777777777777777
$ws.Range("F2").Value = "[1 0 1 1 1 1 1 0 1 1 1 1]"
$ws.Range("H2").Value = 0.6051402927054031
$ws.Range("I2").Value = 0.02169729297097148
$ws.Range("J2").Value = 0.4641835471717825
$ws.Range("K2").Value = 0.07979009093911736

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', RobustScaler()), ('selector', None),`n                ('model', SVC(C=1, kernel='sigmoid', random_state=42))])"
$ws.Range("B3").Value = 0.7522549019607843
$ws.Range("C3").Value = "{'selector': None, 'scaler': RobustScaler(), 'model__kernel': 'sigmoid', 'model__class_weight': None, 'model__C': 1}"
$ws.Range("D3").Value = 0.7058823529411764
$ws.Range("F3").Value = "[1 1 1 1 1 1 0 1 1 1 1 0]"
$ws.Range("H3").Value = 0.6921889557081101
$ws.Range("I3").Value = 0.01881581329958663
$ws.Range("J3").Value = 0.5907850104960398
$ws.Range("K3").Value = 0.06027122592202517

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', StandardScaler()),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model', SVC(C=3, kernel='linear', random_state=42))])"
$ws.Range("B4").Value = 0.7316391941391942
$ws.Range("C4").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': StandardScaler(), 'model__kernel': 'linear', 'model__class_weight': None, 'model__C': 3}"
$ws.Range("H4").Value = 0.6119198291237338
$ws.Range("I4").Value = 0.02245248650787208
$ws.Range("J4").Value = 0.5178457111515935
$ws.Range("K4").Value = 0.07447161774120849

# Row 5
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),`n                ('model',`n                 SVC(C=5, class_weight='balanced', kernel='sigmoid',`n                     random_state=42))])"
$ws.Range("B5").Value = 0.7933006535947712
$ws.Range("C5").Value = "{'selector': None, 'scaler': MinMaxScaler(), 'model__kernel': 'sigmoid', 'model__class_weight': 'balanced', 'model__C': 5}"
$ws.Range("D5").Value = 0.5882352941176471
$ws.Range("F5").Value = "[1 1 1 1 1 1 1 1 1 1 1 1]"
$ws.Range("H5").Value = 0.7720355329864048
$ws.Range("I5").Value = 0.021648373137596
$ws.Range("J5").Value = 0.7027753742989037
$ws.Range("K5").Value = 0.0757205638019951
